$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 758; this pushes the existing
# rows 758-837 down to 760-839 (and extends the used range/dimension
# to A1:R839), matching the new weekly data append.
$ws.Range("A758:A759").EntireRow.Insert()

# --- Populate new row 758 (week of 2023-01-20, $/caja 36 atados) ---
$ws.Cells.Item(758, 1).Value = 9
$ws.Cells.Item(758, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(758, 3).Value = "Metropolitana"
$ws.Cells.Item(758, 4).Value = 44946
$ws.Cells.Item(758, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(758, 5).Value = 13
$ws.Cells.Item(758, 6).Value = 100112040
$ws.Cells.Item(758, 7).Value = "Cilantro"
$ws.Cells.Item(758, 8).Value = "Sin especificar"
$ws.Cells.Item(758, 9).Value = "Primera"
$ws.Cells.Item(758, 10).Value = 70
$ws.Cells.Item(758, 11).Value = 11000
$ws.Cells.Item(758, 12).Value = 11000
$ws.Cells.Item(758, 13).Value = 11000
$ws.Cells.Item(758, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(758, 15).Value = "Región Metropolitana"
$ws.Cells.Item(758, 16).Value = 306
$ws.Cells.Item(758, 17).Value = 36
$ws.Cells.Item(758, 18).Value = "Hortaliza"

# --- Populate new row 759 (week of 2023-01-20, $/docena de atados) ---
$ws.Cells.Item(759, 1).Value = 9
$ws.Cells.Item(759, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(759, 3).Value = "Metropolitana"
$ws.Cells.Item(759, 4).Value = 44946
$ws.Cells.Item(759, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(759, 5).Value = 13
$ws.Cells.Item(759, 6).Value = 100112040
$ws.Cells.Item(759, 7).Value = "Cilantro"
$ws.Cells.Item(759, 8).Value = "Sin especificar"
$ws.Cells.Item(759, 9).Value = "Primera"
$ws.Cells.Item(759, 10).Value = 160
$ws.Cells.Item(759, 11).Value = 24000
$ws.Cells.Item(759, 12).Value = 26000
$ws.Cells.Item(759, 13).Value = 25000
$ws.Cells.Item(759, 14).Value = "$/docena de atados"
$ws.Cells.Item(759, 15).Value = "Región Metropolitana"
$ws.Cells.Item(759, 16).Value = 8333
$ws.Cells.Item(759, 17).Value = 3
$ws.Cells.Item(759, 18).Value = "Hortaliza"
